$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        if ($val -ge 0) {
            $rounded = [math]::Floor($val + 0.5)
        } else {
            $rounded = [math]::Ceiling($val - 0.5)
        }
        $cell.Value = $rounded
    }
}
